# tpivot-phase-3-tracking.xlsx: "ui refactor to match excel idioms"
#
# - Row 24 ("Arrange table fields and sorting buckets in vertical stack") /
#   feature #24: mark Status Complete and stamp Started/Completed dates.
# - Row 25 ("Restyle app for cohesive look and feel") / feature #25: replace
#   the old "Bootstrap would be easiest..." Approach note with the actual
#   writeup of the ribbon-style UI refactor, move Status to "In progress",
#   and stamp a Started date.
# - Move the saved selection/scroll position to reflect where the author was
#   last working (bottom of the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 25 (spreadsheet row 25): Status -> Complete, Started/Completed dates ---
$ws.Range("F25").Value = "Complete"

# Copy number-formatted date style from an existing "date" cell (G3) so the
# new cells land on the same shared cellXf (numFmtId 14) instead of Excel
# inventing a new custom format.
$ws.Range("G3").Copy($ws.Range("G25"))
$ws.Range("G25").Value = 42979

$ws.Range("G3").Copy($ws.Range("H25"))
$ws.Range("H25").Value = 42979

# --- Row 26 (spreadsheet row 26): Approach text, Status -> In progress, Started date ---
$ws.Range("E26").Value = "Redid UI to more closely match Excel ribbon idioms. Used bootstrap for element styles. Still need to peg loading UI to toolbar element."
$ws.Range("F26").Value = "In progress"

$ws.Range("G3").Copy($ws.Range("G26"))
$ws.Range("G26").Value = 42979

# Row grew to a two-line wrap once the longer Approach text was entered.
$ws.Rows.Item(26).RowHeight = 25.5

# --- Restore view/selection to where editing left off ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws.Range("E27").Select()
